$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Medicine" attribute to "First Aid" for the First Aid Kit row
$ws.Range("B6").Value = "First Aid"

# Update the long-form purpose/description text for the First Aid Kit
$ws.Range("D6").Value = "This kit also has bandages and painkillers useful for dealing with less life-threatening scenarios, when coupled with the \imp{Emergency Care} ability, this can be used to remove the \imp{Critical Condition} status. "

# Row height shrinks now that the description is shorter
$ws.Rows.Item(6).RowHeight = 57.45

# Move the active cell selection to D7
$ws.Range("D7").Select()
